$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 303; everything below shifts down by one.
$ws.Rows.Item(303).Insert()

# Populate the newly inserted row 303 with the new data point.
$ws.Range("A303").Value = 10
$ws.Range("B303").Value = "Vega Modelo de Temuco"
$ws.Range("C303").Value = "La Araucanía"
$ws.Range("D303").Value = 45211
$ws.Range("E303").Value = 9
$ws.Range("F303").Value = 100114007
$ws.Range("G303").Value = "Jengibre"
$ws.Range("H303").Value = "Sin especificar"
$ws.Range("I303").Value = "Primera"
$ws.Range("J303").Value = 380
$ws.Range("K303").Value = 22000
$ws.Range("L303").Value = 25000
$ws.Range("M303").Value = 24368
$ws.Range("N303").Value = "`$/caja 13 kilos"
$ws.Range("O303").Value = "Perú"
$ws.Range("P303").Value = 1874
$ws.Range("Q303").Value = 13
$ws.Range("R303").Value = "Hortaliza"
